# Daily COVID-19 Valais data revision: update positive-case, hospitalisation
# and death counts for the last days of the table; dependent totals
# (cumulative cases/deaths, hospitalisation sum) recompute automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns C, E, F, G: plain numeric cells -------------------------------
$ws.Range("E295").Value = 18
$ws.Range("F295").Value = 14
$ws.Range("E296").Value = 17
$ws.Range("F296").Value = 10
$ws.Range("E297").Value = 17
$ws.Range("F297").Value = 12
$ws.Range("E298").Value = 18
$ws.Range("F298").Value = 13
$ws.Range("E299").Value = 18
$ws.Range("F299").Value = 13
$ws.Range("E300").Value = 18
$ws.Range("F300").Value = 15
$ws.Range("E301").Value = 15
$ws.Range("F301").Value = 11
$ws.Range("E302").Value = 15
$ws.Range("F302").Value = 9
$ws.Range("E303").Value = 15
$ws.Range("F303").Value = 12
$ws.Range("E304").Value = 15
$ws.Range("F304").Value = 9
$ws.Range("E305").Value = 15
$ws.Range("F305").Value = 8
$ws.Range("C307").Value = 132
$ws.Range("E307").Value = 15
$ws.Range("F307").Value = 9
$ws.Range("G307").Value = 49
$ws.Range("C308").Value = 118
$ws.Range("C309").Value = 17
$ws.Range("E309").Value = 12
$ws.Range("F309").Value = 8
$ws.Range("G309").Value = 60

# --- Columns L and M are formatted as Text ("@") in this sheet, so writing
# a numeric .Value would be stored as text. Temporarily switch the target
# cells to General, assign the numbers, then restore the Text format so the
# cell style/appearance is unchanged but the stored value is numeric.

# Column L
$ws.Range("L305").NumberFormat = "General"
$ws.Range("L309").NumberFormat = "General"
$ws.Range("L305").Value = 1
$ws.Range("L309").Value = 0
$ws.Range("L305").NumberFormat = "@"
$ws.Range("L309").NumberFormat = "@"

# Column M
$ws.Range("M303").NumberFormat = "General"
$ws.Range("M304").NumberFormat = "General"
$ws.Range("M306").NumberFormat = "General"
$ws.Range("M309").NumberFormat = "General"
$ws.Range("M303").Value = 1
$ws.Range("M304").Value = 2
$ws.Range("M306").Value = 1
$ws.Range("M309").Value = 0
$ws.Range("M303").NumberFormat = "@"
$ws.Range("M304").NumberFormat = "@"
$ws.Range("M306").NumberFormat = "@"
$ws.Range("M309").NumberFormat = "@"

# --- Restore the active cell / selection on the frozen bottom-right pane ---
$ws.Range("T9").Select()
